$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header renames (row 1): 魔法变化 -> MP变化, 领导变化 -> LP变化 ---
$ws.Range("C1").Value = "MP变化"
$ws.Range("D1").Value = "LP变化"

# --- Add new data row 11 (card "巫师学徒" / 57000008) ---
# First clone the formatting of the previous last row (10) down onto row 11
# so the new row picks up the same banded-table styling Excel itself would
# apply when a table is grown by one row.
$ws.Range("A10:I10").Copy()
$ws.Range("A11:I11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = 57000008
$ws.Range("B11").Value = "巫师学徒"
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
# Force these to be stored as literal text ("true"/"false"), matching the
# rest of the boolean-flag columns in this sheet, rather than native
# Excel booleans.
$ws.Range("G11").Value = "'true"
$ws.Range("H11").Value = "'false"
$ws.Range("I11").Value = "'false"

# The leading apostrophe above flips the cell's number format to a
# quote-prefixed text style; re-apply the real formatting from row 10 on
# top (format only) so the new cells end up byte-identical in style to
# their neighbours.
$ws.Range("G10:I10").Copy()
$ws.Range("G11:I11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Grow the table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:I11"))

# --- Selection / view bookkeeping to match the authored state ---
[void]$ws.Range("A11").Select()

try {
    $win = $excel.ActiveWindow
    $win.WindowState = -4140
} catch {
}
